$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 401
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67

# Row 3
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57

# Row 4
$ws.Range("K4").Value = 1.91
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48

# Row 7
$ws.Range("AN7").Value = 5.5
$ws.Range("AW7").Value = 4
$ws.Range("BA7").Value = 67
$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 1.95

# Row 8
$ws.Range("G8").Value = 2.35

# Row 9
$ws.Range("AD9").Value = 5.5
$ws.Range("AW9").Value = 5.5
$ws.Range("K9").Value = 1.91
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 5.5
$ws.Range("O9").Value = 1.52
$ws.Range("P9").Value = 2.37
$ws.Range("S9").Value = 1.6
$ws.Range("T9").Value = 2.25
$ws.Range("U9").Value = 2.15
$ws.Range("V9").Value = 1.62
$ws.Range("X9").Value = 9

# Row 10
$ws.Range("K10").Value = 1.91

# Row 11
$ws.Range("AB11").Value = 29
$ws.Range("AC11").Value = 6.3
$ws.Range("AD11").Value = 6.3
$ws.Range("AE11").Value = 16.5
$ws.Range("AG11").Value = 800
$ws.Range("AI11").Value = 28
$ws.Range("AK11").Value = 100
$ws.Range("AM11").Value = 60
$ws.Range("AN11").Value = 3.55
$ws.Range("AO11").Value = 8.5
$ws.Range("AQ11").Value = 29
$ws.Range("AU11").Value = 7.3
$ws.Range("AV11").Value = 65
$ws.Range("BB11").Value = 450
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 3.25
$ws.Range("J11").Value = 2.27
$ws.Range("K11").Value = 2.1
$ws.Range("N11").Value = 6.3
$ws.Range("O11").Value = 1.37
$ws.Range("P11").Value = 2.85
$ws.Range("U11").Value = 1.93
$ws.Range("V11").Value = 1.78
$ws.Range("X11").Value = 7.8
$ws.Range("Y11").Value = 8
$ws.Range("Z11").Value = 14

# Row 12
$ws.Range("AC12").Value = 5.3
$ws.Range("AD12").Value = 5.1
$ws.Range("AH12").Value = 7.8
$ws.Range("AJ12").Value = 11
$ws.Range("AR12").Value = 100
$ws.Range("AT12").Value = 2.42
$ws.Range("AU12").Value = 6.5
$ws.Range("G12").Value = 2.62
$ws.Range("H12").Value = 2.6
$ws.Range("J12").Value = 3.25
$ws.Range("N12").Value = 5.3
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 2.37
$ws.Range("R12").Value = 1.52
$ws.Range("X12").Value = 12.5

# Row 15
$ws.Range("Q15").Value = 2.3
$ws.Range("R15").Value = 1.6

# Row 16
$ws.Range("K16").Value = 2.38

# Row 17
$ws.Range("AA17").Value = 18
$ws.Range("AB17").Value = 27
$ws.Range("AC17").Value = 9.25
$ws.Range("AD17").Value = 6.1
$ws.Range("AG17").Value = 450
$ws.Range("AJ17").Value = 11
$ws.Range("AL17").Value = 27
$ws.Range("AM17").Value = 32
$ws.Range("AN17").Value = 4.15
$ws.Range("AP17").Value = 18
$ws.Range("AR17").Value = 70
$ws.Range("AS17").Value = 200
$ws.Range("AT17").Value = 2.62
$ws.Range("AU17").Value = 6.4
$ws.Range("AV17").Value = 50
$ws.Range("AX17").Value = 17
$ws.Range("AY17").Value = 21
$ws.Range("AZ17").Value = 80
$ws.Range("BA17").Value = 100
$ws.Range("BB17").Value = 250
$ws.Range("G17").Value = 2.18
$ws.Range("H17").Value = 3.1
$ws.Range("I17").Value = 3.2
$ws.Range("K17").Value = 2.1
$ws.Range("L17").Value = 3.6
$ws.Range("M17").Value = 9.25
$ws.Range("O17").Value = 1.29
$ws.Range("P17").Value = 2.95
$ws.Range("Q17").Value = 1.93
$ws.Range("R17").Value = 1.78
$ws.Range("W17").Value = 7.5
$ws.Range("Y17").Value = 8.75
$ws.Range("Z17").Value = 22

# Row 18
$ws.Range("G18").Value = 1.6
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13
$ws.Range("Q18").Value = 1.83
$ws.Range("R18").Value = 2.03

# Row 19
$ws.Range("AE19").Value = 21
$ws.Range("AL19").Value = 41
$ws.Range("AS19").Value = 126
$ws.Range("AT19").Value = 3.25
$ws.Range("AU19").Value = 9
$ws.Range("AV19").Value = 51
$ws.Range("G19").Value = 1.44
$ws.Range("L19").Value = 6
$ws.Range("O19").Value = 1.22
$ws.Range("P19").Value = 4
$ws.Range("Q19").Value = 1.73
$ws.Range("R19").Value = 2.08
$ws.Range("S19").Value = 1.33
$ws.Range("T19").Value = 3.25
$ws.Range("U19").Value = 2
$ws.Range("V19").Value = 1.73
$ws.Range("W19").Value = 7
$ws.Range("X19").Value = 7

# Row 20
$ws.Range("I20").Value = 1.9
$ws.Range("V20").Value = 1.73

# Row 21
$ws.Range("AA21").Value = 41
$ws.Range("AK21").Value = 10
$ws.Range("AN21").Value = 8
$ws.Range("AX21").Value = 7
$ws.Range("G21").Value = 6.5
$ws.Range("I21").Value = 1.42
$ws.Range("M21").Value = 1.01
$ws.Range("N21").Value = 15
$ws.Range("U21").Value = 1.8
$ws.Range("V21").Value = 1.91
$ws.Range("X21").Value = 34
$ws.Range("Z21").Value = 67
